# SectorGroup.xlsx fix: the codeforiati:group-code / group-name / category-name
# columns (E, F, G) were shifted by one position. For every row (including the
# header), the correct values are:
#   new E = old G
#   new F = old E
#   new G = old F
# (a right-rotation of the E,F,G triple). Columns A-D are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$colE = 5
$colF = 6
$colG = 7

# First pass: snapshot the current E, F, G values for every row so the
# write-back (second pass) never reads an already-updated cell.
$origE = @{}
$origF = @{}
$origG = @{}

for ($r = 1; $r -le $lastRow; $r++) {
    $origE[$r] = $ws.Cells.Item($r, $colE).Value2
    $origF[$r] = $ws.Cells.Item($r, $colF).Value2
    $origG[$r] = $ws.Cells.Item($r, $colG).Value2
}

# Second pass: write the rotated values back.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colE).Value = $origG[$r]
    $ws.Cells.Item($r, $colF).Value = $origE[$r]
    $ws.Cells.Item($r, $colG).Value = $origF[$r]
}

Write-Output "Rotated E/F/G for $lastRow rows"
